# Commit message: "Fruta / hortaliza, semanal" — weekly refresh of the
# "Ají" price sheet: two brand-new price entries (one "Primera", one
# "Segunda" quality, both dated 2023-10-10 / serial 45209, Región de
# Arica y Parinacota origin) are inserted at the top of the Terminal La
# Palmera de La Serena - Ají data block (row 467 onward), pushing all
# the existing data rows down by two and growing the used range from
# A1:R564 to A1:R566.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 467-468; everything currently on/after
# row 467 (old rows 467..564) shifts down to 469..566.
$ws.Rows("467:468").Insert()

# ---- New row 467 : Inferno / Primera --------------------------------
$ws.Cells.Item(467, 1).Value  = 8
$ws.Cells.Item(467, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(467, 3).Value  = "Coquimbo"
$ws.Cells.Item(467, 4).Value  = 45209
$ws.Cells.Item(467, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(467, 5).Value  = 4
$ws.Cells.Item(467, 6).Value  = 100112021
$ws.Cells.Item(467, 7).Value  = "Ají"
$ws.Cells.Item(467, 8).Value  = "Inferno"
$ws.Cells.Item(467, 9).Value  = "Primera"
$ws.Cells.Item(467, 10).Value = 400
$ws.Cells.Item(467, 11).Value = 29000
$ws.Cells.Item(467, 12).Value = 30000
$ws.Cells.Item(467, 13).Value = 29500
$ws.Cells.Item(467, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(467, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(467, 16).Value = 2950
$ws.Cells.Item(467, 17).Value = 10
$ws.Cells.Item(467, 18).Value = "Hortaliza"

# ---- New row 468 : Inferno / Segunda ---------------------------------
$ws.Cells.Item(468, 1).Value  = 8
$ws.Cells.Item(468, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(468, 3).Value  = "Coquimbo"
$ws.Cells.Item(468, 4).Value  = 45209
$ws.Cells.Item(468, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(468, 5).Value  = 4
$ws.Cells.Item(468, 6).Value  = 100112021
$ws.Cells.Item(468, 7).Value  = "Ají"
$ws.Cells.Item(468, 8).Value  = "Inferno"
$ws.Cells.Item(468, 9).Value  = "Segunda"
$ws.Cells.Item(468, 10).Value = 240
$ws.Cells.Item(468, 11).Value = 19000
$ws.Cells.Item(468, 12).Value = 20000
$ws.Cells.Item(468, 13).Value = 19500
$ws.Cells.Item(468, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(468, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(468, 16).Value = 1950
$ws.Cells.Item(468, 17).Value = 10
$ws.Cells.Item(468, 18).Value = "Hortaliza"
